$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width updates ---
# Note: Excel COM's ColumnWidth property (Calibri-11 character units) gets
# re-quantized to pixels internally, which shifts the raw XML 'width'
# attribute by ~0.83 relative to the character value you assign. The
# offsets below were calibrated so the saved XML 'width' lands exactly on
# the target integer values (43 / 39 / 16 / 15 / 47).
$ws.Columns.Item(3).ColumnWidth = 42.17
$ws.Columns.Item(4).ColumnWidth = 38.17
$ws.Columns.Item(6).ColumnWidth = 15.17
$ws.Columns.Item(7).ColumnWidth = 14.17
$ws.Columns.Item(8).ColumnWidth = 46.17

# Column A holds numeric-looking opportunity IDs that must stay stored as
# text (not auto-converted to numbers) -- force text format before writing
# so Excel keeps them as strings, then restore the default style so no
# stray formatting is left behind.
$ws.Range("A2:A5").NumberFormat = "@"

# --- Row 2 ---
$ws.Range("A2").Value = "1326927"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1326927"
$ws.Range("C2").Value = "Marketing Expert"
$ws.Range("D2").Value = "Visakhapatnam, Andhra Pradesh, India"
$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "0 applicants"
$ws.Range("G2").Value = "9 - 12 Weeks"
$ws.Range("H2").Value = "Think Big Global Education Solution Pvt Ltd."

# --- Row 3 ---
$ws.Range("A3").Value = "1326926"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1326926"
$ws.Range("C3").Value = "Study Abroad Counsellor"
$ws.Range("D3").Value = "Visakhapatnam, Andhra Pradesh, India"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "0 applicants"
$ws.Range("G3").Value = "9 - 12 Weeks"
$ws.Range("H3").Value = "Think Big Global Education Solution Pvt Ltd."

# --- Row 4 ---
$ws.Range("A4").Value = "1326925"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1326925"
$ws.Range("C4").Value = "Digital Marketer"
$ws.Range("D4").Value = "Visakhapatnam, Andhra Pradesh, India"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "0 applicants"
$ws.Range("G4").Value = "9 - 12 Weeks"
$ws.Range("H4").Value = "Think Big Global Education Solution Pvt Ltd."

# --- Row 5 ---
$ws.Range("A5").Value = "1325700"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1325700"
$ws.Range("C5").Value = "Business Analyst and Executive Secretary"
$ws.Range("D5").Value = "Colombo, Sri Lanka"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "20 applicants"
$ws.Range("G5").Value = "3 - 6 Months"
$ws.Range("H5").Value = "Indian Kitchen PVT LTD"

# Restore the default (unstyled) cell style for column A now that the
# text-formatted values have been written.
$ws.Range("A2:A5").Style = "Normal"

# --- Row 6 no longer exists in the source data: delete it so the sheet
#     dimension shrinks back down to row 5 ---
$ws.Rows.Item(6).Delete()
